# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables "header" rows embedded as plain text in the first
# row(s) of each worksheet use PascalCase attribute names
# (Type=, Id=, Name=, Description=, Date=, ObjTablesVersion=, TableID=,
# TableName=). They must be rewritten to lowerCamelCase
# (type=, id=, name=, description=, date=, objTablesVersion=, tableID=,
# tableName=).

$wb = $excel.ActiveWorkbook

# Sheet "!!_Table of contents" -> row 1 (document-level marker) and
# row 2 (table-level marker).
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' tableID='Table of contents' tableName='Readme' description='Table/model and column/attribute definitions' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# Sheet "!!Child" -> row 1 (table-level marker).
$wsChild = $wb.Worksheets.Item("!!Child")
$wsChild.Range("A1").Value = "!!ObjTables type='Data' id='Child' name='Child' description='Represents a child' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"

# Sheet "!!Parent" -> row 1 (table-level marker).
$wsParent = $wb.Worksheets.Item("!!Parent")
$wsParent.Range("A1").Value = "!!ObjTables type='Data' id='Parent' name='Parent' description='Represents a parent' date='2019-09-18 00:10:05' objTablesVersion='0.0.8'"
